$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# The shape's text body ends with an empty paragraph. Appending to the
# whole-shape TextRange lands the new text inside that trailing empty
# paragraph, right before its existing endParaRPr - exactly where the
# new sentence belongs.
$startLen = $tr.Text.Length

$run1 = $tr.InsertAfter("The number of ")
$run2 = $tr.InsertAfter("guess is set to 6.")

Write-Host "Full text now: $($tr.Text)"
